# Adds a new data-collection wave (26. 1. 2021) as an extra trailing column
# on both worksheets ("data" and "pocetR"), and bumps the "aktualizace"
# date in the footer label cells from 12. 1. 2021 to 2. 2. 2021.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("pocetR")

# ---------------------------------------------------------------------
# Sheet "data": new column X (percentages), header + 22 data rows
# ---------------------------------------------------------------------

# Header cell: copy style from the previous header cell (W1) so the new
# header keeps the bold/centered/bordered look, then set its text.
$ws1.Range("W1").Copy()
$ws1.Range("X1").PasteSpecial(-4122)
$ws1.Range("X1").Value = "26. 1. 2021"

$sheet1Values = @{
    2  = 0.12
    3  = 0.09
    4  = 0.14
    5  = 0.07000000000000001
    6  = 0.1
    7  = 0.14
    8  = 0.14
    9  = 0.21
    10 = 0.13
    11 = 0.09
    12 = 0.13
    13 = 0.1
    14 = 0.25
    15 = 0.13
    16 = 0.09
    17 = 0.18
    18 = 0.13
    19 = 0.09
    20 = 0.07000000000000001
    21 = 0.07000000000000001
    22 = 0.07000000000000001
    23 = 0.17
}

foreach ($row in $sheet1Values.Keys) {
    $ws1.Cells.Item($row, 24).Value = $sheet1Values[$row]
}

# Footer label (row 24, column A): bump the "aktualizace" date.
$ws1.Range("A24").Value = "Život během pandemie, Duševní zdraví, % respondentů celkově a ve skupinách, aktualizace 2. 2. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR": new column W (sample sizes), header + 22 data rows
# ---------------------------------------------------------------------

$ws2.Range("V1").Copy()
$ws2.Range("W1").PasteSpecial(-4122)
$ws2.Range("W1").Value = "26. 1. 2021"

$sheet2Values = @{
    2  = 2131
    3  = 1040
    4  = 1091
    5  = 287
    6  = 732
    7  = 319
    8  = 764
    9  = 173
    10 = 329
    11 = 402
    12 = 379
    13 = 848
    14 = 237
    15 = 491
    16 = 1403
    17 = 272
    18 = 829
    19 = 675
    20 = 203
    21 = 388
    22 = 774
    23 = 969
}

foreach ($row in $sheet2Values.Keys) {
    $ws2.Cells.Item($row, 23).Value = $sheet2Values[$row]
}

# Row 24 on this sheet is a trailing "blank" row (string-typed empty
# cells across B:V) under the footer label in column A. Extend that same
# blank-cell pattern into the new column W to match the existing row.
$ws2.Range("V24").Copy()
$ws2.Range("W24").PasteSpecial(-4104)

# Footer label (row 24, column A): bump the "aktualizace" date.
$ws2.Range("A24").Value = "Život během pandemie, Duševní zdraví, velikost dotázaného souboru celkově a ve skupinách, aktualizace 2. 2. 2021"
